$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates per diff (cryptos.xlsx symbol-list refresh).
# NumberFormat is forced to text ("@") before writing so numeric-looking
# strings (prices, volume codes) are preserved as text, matching the
# original inlineStr cell type rather than being coerced to a float.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '244.26'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '23.81'
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = 'LEO'
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '3.555'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '3LEOLEO'
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = 'HuobiToken'
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '5.304'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '4HuobiTokenHT'
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = 'Cronos'
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.05864'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '5CronosCRO'
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '6.476'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '6KuCoinTokenKCS'
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = 'GateToken'
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.342'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '7GateTokenGT'
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8162'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '8MXTokenMX'
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = 'FTXToken'
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8925'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '9FTXTokenFTT'
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1388'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '10WazirXWRX'
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07254'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '11MandalaExchangeTokenMDX'
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03100'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '12LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.03021'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '13BitrueCoinBTR'
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.09361'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '14BitMartTokenBMX'
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'MCDex'
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.833'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '15MCDexMCB'
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'BitForexToken'
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.001551'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '16BitForexTokenBF'
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = 'CoinExToken'
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.04729'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '17CoinExTokenCET'
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = 'One'
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0006038'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '18OneONE'
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = 'TigerCash'
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.006265'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '19TigerCashTCH'
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = 'BitKan'
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.001264'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '20BitKanKAN'
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = 'HotbitToken'
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.004595'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '21HotbitTokenHTB'
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = 'NitroEx'
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.00008727'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '22NitroExNTX'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.176'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1313'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0002347'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006310'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1057'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002550'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007128'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005330'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000752'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5417'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '46CoinbaseStockTokenCOINWorstin24h'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.01837'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002106'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0002006'
